$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New set of indicator rows (labels + values) replacing the old rows 2-7
$rows = @(
    @("Mean Square Error", 2243454.062),
    @("Root Mean Square Error", 1497.816),
    @("Mean Absolute Error", 1272.794),
    @("Root Mean Square Error (log)", 0.047),
    @("Mean Absolute Error (log)", 0.033),
    @("Mean Absolute Percentage Error", 6.272),
    @("U-Theil1 Statistic", 67.66800000000001),
    @("U-Theil2 Statistic", 0.21),
    @("Diebold-Mariano Test", 35.591),
    @("Rendement Absolue", 0.35)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
